$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 310 - this shifts rows 310:367 down to 311:368
$ws.Rows.Item(310).EntireRow.Insert()

# Populate the newly inserted row 310 with the new record's data
$ws.Range("A310").Value = 3
$ws.Range("B310").Value = "Femacal de La Calera"
$ws.Range("C310").Value = "Coquimbo"
$ws.Range("D310").Value = 44694
$ws.Range("E310").Value = 5
$ws.Range("F310").Value = 100112031
$ws.Range("G310").Value = "Poroto verde"
$ws.Range("H310").Value = "Magnum"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 73
$ws.Range("K310").Value = 28000
$ws.Range("L310").Value = 29000
$ws.Range("M310").Value = 28521
$ws.Range("N310").Value = "`$/malla 25 kilos"
$ws.Range("O310").Value = "Provincia de Santiago"
$ws.Range("P310").Value = 1141
$ws.Range("Q310").Value = 25
$ws.Range("R310").Value = "Hortaliza"
